$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.462512195110321
$ws.Range("B1").Value = 1.316810131072998
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.060950756072998
$ws.Range("E1").Value = 1.124637722969055
